$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill marker info ("NAT") into the marker_1 column (K) for status 8 rows
# that are currently missing it.
$ws.Range("K13").Value = "NAT"
$ws.Range("K14").Value = "NAT"
$ws.Range("K15").Value = "NAT"
$ws.Range("K18").Value = "NAT"
$ws.Range("K19").Value = "NAT"

# Leave the selection on the last filled cell, matching the saved view state.
$ws.Range("K18").Select()
